# Update automático del index.html y archivo Excel
# Inserta un nuevo reclamo (caso 5651) encima de la fila del caso 5887,
# desplazando hacia abajo todas las filas siguientes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserta una fila nueva en la posición 27 (empuja 27..56 -> 28..57)
$ws.Rows(27).Insert()

# --- Completa los datos de la nueva fila 27 ---

# Columnas que contienen texto "numérico" o con forma de fecha: hay que
# forzarlas a texto (prefijo con apóstrofe) y luego limpiar el formato que
# ese prefijo deja aplicado, para que el valor final quede como texto plano
# sin ningún formato de número especial.
$ws.Range("A27").Value = "'5651"
$ws.Range("A27").ClearFormats()

$ws.Range("B27").Value = "'4/22/2025"
$ws.Range("B27").ClearFormats()

$ws.Range("D27").Value = "'4"
$ws.Range("D27").ClearFormats()

$ws.Range("E27").Value = "'804876051"
$ws.Range("E27").ClearFormats()

$ws.Range("I27").Value = "'1"
$ws.Range("I27").ClearFormats()

# Columnas de texto normal (no requieren protección)
$ws.Range("C27").Value = "MONTES DE OCA, MANUEL AV. 511"
$ws.Range("F27").Value = "Optical Power"
$ws.Range("G27").Value = "Pendiente"
$ws.Range("H27").Value = "Pegar los ductos al prfv"
$ws.Range("J27").Value = "Cambio"
$ws.Range("K27").Value = "Sin equipos"
$ws.Range("L27").Value = "Pasante"

# Coordenadas (numéricas)
$ws.Range("M27").Value = -58.375515
$ws.Range("N27").Value = -34.634393
